# Financials update: insert a new first data column (D) before the existing
# quarterly/yearly columns and populate it with the latest period's figures.
# All previously existing columns D:K shift right to E:L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column before column D; existing D:K -> E:L.
$ws.Columns("D:D").Insert()

# The freshly inserted column has no formatting yet (it picks up the format
# of the column to its left). Copy number formats from the (now-shifted)
# original column E so the new column D matches the rest of the data block
# (date format on the header rows, thousands-separator number format on the
# data rows).
$ws.Range("E5:K102").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match column width to its neighboring data columns.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 3654100
$ws.Range("D9").Value = 3129400
$ws.Range("D10").Value = 524700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -52400
$ws.Range("D15").Value = 168700
$ws.Range("D17").Value = 3488200
$ws.Range("D18").Value = 165900
$ws.Range("D20").Value = 400
$ws.Range("D21").Value = 335000
$ws.Range("D22").Value = 99200
$ws.Range("D23").Value = 67100
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 67000
$ws.Range("D27").Value = -9300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -400
$ws.Range("D33").Value = -9300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -9300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 17200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 251500
$ws.Range("D44").Value = 64600
$ws.Range("D45").Value = 46000
$ws.Range("D46").Value = 379300
$ws.Range("D47").Value = 1188200
$ws.Range("D48").Value = 2029700
$ws.Range("D49").Value = 692400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 4900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4294500
$ws.Range("D57").Value = 213000
$ws.Range("D58").Value = 3300
$ws.Range("D59").Value = 115800
$ws.Range("D60").Value = 332100
$ws.Range("D61").Value = 1752400
$ws.Range("D62").Value = 176200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2442000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 612000
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1240500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -9300
$ws.Range("D83").Value = 168700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 253600
$ws.Range("D91").Value = -305500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -241200
$ws.Range("D96").Value = -230900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 3500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 15900
